# Alumni route refactor: update employment status breakdown by industry,
# fix names, add per-row "company_name -> NA" fallback, and restyle
# the enrollment_date column header / row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update data rows (columns A, C, D, E). Column B (enrollment_date)
#    values are unchanged row for row, so we leave them untouched.
# ---------------------------------------------------------------------

$data = @(
    @{ Row = 2;  Name = "Raman";   Status = "Employed";   Industry = "Healthcare"; Company = "FinBank" },
    @{ Row = 3;  Name = "Shayam";  Status = "Employed";   Industry = "IT";         Company = "FinBank" },
    @{ Row = 4;  Name = "Karan";   Status = "UnEmployed"; Industry = "IT";         Company = "NA" },
    @{ Row = 5;  Name = "rasat";   Status = "UnEmployed"; Industry = "Healthcare"; Company = "NA" },
    @{ Row = 6;  Name = "sdfa";    Status = "UnEmployed"; Industry = "Education";  Company = "NA" },
    @{ Row = 7;  Name = "dsfs";    Status = "UnEmployed"; Industry = "Finance";    Company = ("NA" + [char]10 + "NA" + [char]10 + "NA") },
    @{ Row = 8;  Name = "sdfsd";   Status = "UnEmployed"; Industry = "Finance";    Company = "NA" },
    @{ Row = 9;  Name = "fsfd";    Status = "UnEmployed"; Industry = "IT";         Company = "NA" },
    @{ Row = 10; Name = "sdfsd";   Status = "Employed";   Industry = "IT";         Company = "ShopEase" },
    @{ Row = 11; Name = "sdfsdf";  Status = "UnEmployed"; Industry = "Finance";    Company = "NA" }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.Name
    $ws.Cells.Item($r, 3).Value = $rec.Status
    $ws.Cells.Item($r, 4).Value = $rec.Industry
    $ws.Cells.Item($r, 5).Value = $rec.Company
}

# ---------------------------------------------------------------------
# 2. Header restyle: enrollment_date header (B1) becomes left aligned.
# ---------------------------------------------------------------------

$ws.Range("B1").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 3. Row heights for the data rows grow slightly (18.75 -> 19.5).
# ---------------------------------------------------------------------

for ($r = 2; $r -le 11; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}

# ---------------------------------------------------------------------
# 4. Explicit (automatic/black) font color on the enrollment_date data
#    column, matching the font recolor captured in the workbook theme.
# ---------------------------------------------------------------------

$ws.Range("B2:B11").Font.Color = 0

$wb.Save()
